$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1749.1666
$ws.Range("I4").Value = 1749.1666
$ws.Range("K4").Value = 1749.1666
$ws.Range("M4").Value = -1635.1666
$ws.Range("H10").Value = 5407.3335
$ws.Range("I10").Value = 3111
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 3111
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -2818
$ws.Range("N10").Value = -10586
$ws.Range("H19").Value = 1047.1765
$ws.Range("I19").Value = 425
$ws.Range("K19").Value = 425
$ws.Range("M19").Value = -250
$ws.Range("H28").Value = 1002.1053
$ws.Range("I28").Value = 1114.6428
$ws.Range("J28").Value = 687
$ws.Range("K28").Value = 1114.6428
$ws.Range("L28").Value = 687
$ws.Range("M28").Value = -629.6428000000001
$ws.Range("N28").Value = -1657
$ws.Range("H43").Value = 8502.333000000001
$ws.Range("I43").Value = 3276.2856
$ws.Range("J43").Value = 13075.125
$ws.Range("K43").Value = 3276.2856
$ws.Range("L43").Value = 13075.125
$ws.Range("M43").Value = -3207.2856
$ws.Range("N43").Value = -13213.125
$ws.Range("H138").Value = 8930017
$ws.Range("I138").Value = 1347.2609
$ws.Range("J138").Value = 50001900
$ws.Range("K138").Value = 4041.7827
$ws.Range("L138").Value = 150005700
$ws.Range("M138").Value = 1098.2173
$ws.Range("N138").Value = -150015980

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = -432
$ws.Range("H5").Value = 1035
$ws.Range("I5").Value = 70.5
$ws.Range("J5").Value = 1999.5
$ws.Range("K5").Value = 70.5
$ws.Range("L5").Value = 1999.5
$ws.Range("M5").Value = 41.5
$ws.Range("N5").Value = -2223.5
$ws.Range("H32").Value = 7367.1772
$ws.Range("I32").Value = 7863.635
$ws.Range("K32").Value = 7863.635
$ws.Range("M32").Value = -7576.635

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1035
$ws.Range("I4").Value = 70.5
$ws.Range("J4").Value = 1999.5
$ws.Range("K4").Value = 70.5
$ws.Range("L4").Value = 1999.5
$ws.Range("M4").Value = 44.5
$ws.Range("N4").Value = -2229.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 88500340
$ws.Range("I4").Value = 166667330
$ws.Range("K4").Value = 166667330
$ws.Range("M4").Value = -166667218
$ws.Range("H7").Value = 138.85715
$ws.Range("I7").Value = 104.4
$ws.Range("J7").Value = 225
$ws.Range("K7").Value = 104.4
$ws.Range("L7").Value = 225
$ws.Range("M7").Value = 8.599999999999994
$ws.Range("N7").Value = -451
$ws.Range("H31").Value = 10106574
$ws.Range("I31").Value = 5964.4614
$ws.Range("J31").Value = 47623124
$ws.Range("K31").Value = 5964.4614
$ws.Range("L31").Value = 47623124
$ws.Range("M31").Value = -5669.4614
$ws.Range("N31").Value = -47623714
$ws.Range("H34").Value = 10106574
$ws.Range("I34").Value = 5964.4614
$ws.Range("J34").Value = 47623124
$ws.Range("K34").Value = 5964.4614
$ws.Range("L34").Value = 47623124
$ws.Range("M34").Value = -5762.4614
$ws.Range("N34").Value = -47623528
$ws.Range("H62").Value = 2249.4443
$ws.Range("I62").Value = 2293.75
$ws.Range("J62").Value = 1895
$ws.Range("K62").Value = 2293.75
$ws.Range("L62").Value = 1895
$ws.Range("M62").Value = -1669.75
$ws.Range("N62").Value = -3143
$ws.Range("H65").Value = 2249.4443
$ws.Range("I65").Value = 2293.75
$ws.Range("J65").Value = 1895
$ws.Range("K65").Value = 11468.75
$ws.Range("L65").Value = 9475
$ws.Range("M65").Value = -8348.75
$ws.Range("N65").Value = -15715
$ws.Range("H134").Value = 1870.129
$ws.Range("I134").Value = 1665.4286
$ws.Range("J134").Value = 2300
$ws.Range("K134").Value = 4996.2858
$ws.Range("L134").Value = 6900
$ws.Range("M134").Value = -2461.2858
$ws.Range("N134").Value = -11970

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 764.95
$ws.Range("I4").Value = 183.33333
$ws.Range("J4").Value = 5999.5
$ws.Range("K4").Value = 549.99999
$ws.Range("L4").Value = 17998.5
$ws.Range("M4").Value = -437.99999
$ws.Range("N4").Value = -18222.5
$ws.Range("H131").Value = 1356355.6
$ws.Range("I131").Value = 10392287
$ws.Range("J131").Value = 965.875
$ws.Range("K131").Value = 31176861
$ws.Range("L131").Value = 2897.625
$ws.Range("M131").Value = -31171821
$ws.Range("N131").Value = -12977.625
$ws.Range("H140").Value = 5175.4116
$ws.Range("I140").Value = 2468.75
$ws.Range("K140").Value = 7406.25
$ws.Range("M140").Value = -2226.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3004
$ws.Range("I5").Value = 3004
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3004
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -2892

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5002625
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 5002625
$ws.Range("K2").Value = 0
$ws.Range("L2").ClearContents()
$ws.Range("M2").Value = 5002625
$ws.Range("N2").Value = -5002849
$ws.Range("H7").Value = 11125
$ws.Range("I7").Value = 15375
$ws.Range("J7").Value = 6875
$ws.Range("K7").Value = 15375
$ws.Range("L7").Value = 6875
$ws.Range("M7").Value = -15263
$ws.Range("N7").Value = -7099
$ws.Range("H68").Value = 1915.1666
$ws.Range("I68").Value = 1230.5
$ws.Range("J68").Value = 2257.5
$ws.Range("K68").Value = 1230.5
$ws.Range("L68").Value = 2257.5
$ws.Range("M68").Value = -481.5
$ws.Range("N68").Value = -3755.5
$ws.Range("H71").Value = 1915.1666
$ws.Range("I71").Value = 1230.5
$ws.Range("J71").Value = 2257.5
$ws.Range("K71").Value = 6152.5
$ws.Range("L71").Value = 11287.5
$ws.Range("M71").Value = -2408.5
$ws.Range("N71").Value = -18775.5
$ws.Range("H126").Value = 11125
$ws.Range("I126").Value = 15375
$ws.Range("J126").Value = 6875
$ws.Range("K126").Value = 46125
$ws.Range("L126").Value = 20625
$ws.Range("M126").Value = -43655
$ws.Range("N126").Value = -25565

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 555
$ws.Range("J2").Value = 555
$ws.Range("L2").Value = 555
$ws.Range("N2").Value = -779
$ws.Range("H46").Value = 39655.8
$ws.Range("J46").Value = 39655.8
$ws.Range("L46").Value = 39655.8
$ws.Range("N46").Value = -40117.8
$ws.Range("H132").Value = 2410.1667
$ws.Range("I132").Value = 1827.5
$ws.Range("K132").Value = 5482.5
$ws.Range("M132").Value = -2952.5
$ws.Range("H134").Value = 39655.8
$ws.Range("J134").Value = 39655.8
$ws.Range("L134").Value = 118967.4
$ws.Range("N134").Value = -124037.4
